$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full table data (rows 2-18), reflecting roster reshuffle:
# a new player (Stephon Castle) was added and "Tim Hardaway Jr." removed,
# with the remaining rows re-sorted.
$data = @(
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Lonzo Ball", "PG", "Chicago Bulls"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
